$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New daily rows (15, 16, 17 aprile 2021) appended after row 226
$newRows = @(
    @{ Row = 227; Date = 44301; B = 0; C = 8; D = 326.3973888208894 },
    @{ Row = 228; Date = 44302; B = 0; C = 7; D = 285.5977152182783 },
    @{ Row = 229; Date = 44303; B = 0; C = 4; D = 163.1986944104447 }
)

foreach ($r in $newRows) {
    $rowNum = $r.Row

    $cellA = $ws.Cells.Item($rowNum, 1)
    $cellA.Value = $r.Date
    $ws.Cells.Item($rowNum - 1, 1).Copy()
    $cellA.PasteSpecial(-4122)  # xlPasteFormats

    $ws.Cells.Item($rowNum, 2).Value = $r.B
    $ws.Cells.Item($rowNum, 3).Value = $r.C
    $ws.Cells.Item($rowNum, 4).Value = $r.D
}
